$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tue Jan 28 22:08:22 EST 2025"
$ws.Range("B3").Value = "Tue Jan 28 22:08:36 EST 2025"
$ws.Range("B4").Value = "Tue Jan 28 22:08:48 EST 2025"
$ws.Range("B5").Value = "Tue Jan 28 22:09:02 EST 2025"
$ws.Range("B6").Value = "Tue Jan 28 22:09:15 EST 2025"
$ws.Range("B7").Value = "Tue Jan 28 22:09:28 EST 2025"

$ws.Range("A6").Value = "Pass"
$ws.Range("A7").Value = "Pass"

$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

$ws.Range("C7").Select()
